# Update the cryptocurrency price/volume snapshot (GitHub Actions scheduled refresh).
# Column D holds the "Price" text and column E holds the "Volume(1h)" percentage text.
# These are stored as text (not numbers) in the workbook, so we force the NumberFormat
# of the Price column to "@" (Text) before writing, which prevents Excel from
# auto-converting numeric-looking strings (e.g. "9.000", "0.02360", "24.248.62")
# into floating point numbers and losing their exact original formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2,4).Value  = "24.248.62"
$ws.Cells.Item(2,5).Value  = "  +11.68%  "

$ws.Cells.Item(3,4).Value  = "1.683.65"
$ws.Cells.Item(3,5).Value  = "  +6.94%  "

$ws.Cells.Item(4,4).Value  = "1.003"
$ws.Cells.Item(4,5).Value  = "  +0.50%  "

$ws.Cells.Item(5,4).Value  = "309.64"
$ws.Cells.Item(5,5).Value  = "  +8.89%  "

$ws.Cells.Item(6,4).Value  = "0.9977"
$ws.Cells.Item(6,5).Value  = "  +2.34%  "

$ws.Cells.Item(7,4).Value  = "0.3751"
$ws.Cells.Item(7,5).Value  = "  +1.89%  "

$ws.Cells.Item(8,4).Value  = "0.3457"
$ws.Cells.Item(8,5).Value  = "  +5.93%  "

$ws.Cells.Item(9,4).Value  = "47.83"
$ws.Cells.Item(9,5).Value  = "  +16.10%  "

$ws.Cells.Item(10,4).Value = "1.197"
$ws.Cells.Item(10,5).Value = "  +5.68%  "

$ws.Cells.Item(11,4).Value = "0.07329"
$ws.Cells.Item(11,5).Value = "  +3.99%  "

$ws.Cells.Item(12,4).Value = "0.9993"
$ws.Cells.Item(12,5).Value = "  +0.50%  "

$ws.Cells.Item(13,4).Value = "20.59"
$ws.Cells.Item(13,5).Value = "  +2.40%  "

$ws.Cells.Item(14,4).Value = "6.138"
$ws.Cells.Item(14,5).Value = "  +5.98%  "

$ws.Cells.Item(15,4).Value = "6.801"
$ws.Cells.Item(15,5).Value = "  +5.00%  "

$ws.Cells.Item(16,4).Value = "1.680.25"
$ws.Cells.Item(16,5).Value = "  +7.47%  "

$ws.Cells.Item(17,4).Value = "0.00001114"
$ws.Cells.Item(17,5).Value = "  +4.19%  "

$ws.Cells.Item(18,4).Value = "0.9974"
$ws.Cells.Item(18,5).Value = "  +2.53%  "

$ws.Cells.Item(19,4).Value = "0.06731"
$ws.Cells.Item(19,5).Value = "  +8.88%  "

$ws.Cells.Item(20,4).Value = "82.38"
$ws.Cells.Item(20,5).Value = "  +11.52%  "

$ws.Cells.Item(21,4).Value = "16.57"
$ws.Cells.Item(21,5).Value = "  +3.97%  "

$ws.Cells.Item(22,4).Value = "6.137"
$ws.Cells.Item(22,5).Value = "  +5.25%  "

$ws.Cells.Item(23,5).Value = "  +4.28%  "

$ws.Cells.Item(24,4).Value = "24.213.41"
$ws.Cells.Item(24,5).Value = "  +11.55%  "

$ws.Cells.Item(25,4).Value = "2.411"
$ws.Cells.Item(25,5).Value = "  +3.72%  "

$ws.Cells.Item(26,4).Value = "2.695"
$ws.Cells.Item(26,5).Value = "  +12.23%  "

$ws.Cells.Item(27,4).Value = "3.364"
$ws.Cells.Item(27,5).Value = "  -8.84%  "

$ws.Cells.Item(28,4).Value = "152.49"
$ws.Cells.Item(28,5).Value = "  +2.66%  "

$ws.Cells.Item(29,4).Value = "19.67"
$ws.Cells.Item(29,5).Value = "  +8.21%  "

$ws.Cells.Item(30,4).Value = "1.864.22"
$ws.Cells.Item(30,5).Value = "  +7.39%  "

$ws.Cells.Item(31,4).Value = "127.29"
$ws.Cells.Item(31,5).Value = "  +6.03%  "

$ws.Cells.Item(32,4).Value = "6.497"
$ws.Cells.Item(32,5).Value = "  +20.76%  "

$ws.Cells.Item(33,4).Value = "4.068"
$ws.Cells.Item(33,5).Value = "  +0.34%  "

$ws.Cells.Item(34,4).Value = "0.9966"
$ws.Cells.Item(34,5).Value = "  +10.25%  "

$ws.Cells.Item(35,4).Value = "1.794"
$ws.Cells.Item(35,5).Value = "  +14.30%  "

$ws.Cells.Item(36,4).Value = "0.08509"
$ws.Cells.Item(36,5).Value = "  +4.28%  "

$ws.Cells.Item(37,4).Value = "12.63"
$ws.Cells.Item(37,5).Value = "  +9.68%  "

$ws.Cells.Item(38,4).Value = "0.06512"
$ws.Cells.Item(38,5).Value = "  +8.33%  "

$ws.Cells.Item(39,4).Value = "5.416"
$ws.Cells.Item(39,5).Value = "  +6.76%  "

$ws.Cells.Item(40,4).Value = "9.000"
$ws.Cells.Item(40,5).Value = "  +10.79%  "

$ws.Cells.Item(41,4).Value = "0.02360"
$ws.Cells.Item(41,5).Value = "  +9.24%  "

$ws.Cells.Item(42,4).Value = "1.288"
$ws.Cells.Item(42,5).Value = "  +4.29%  "

$ws.Cells.Item(43,4).Value = "0.2150"
$ws.Cells.Item(43,5).Value = "  +7.18%  "

$ws.Cells.Item(44,4).Value = "0.6231"
$ws.Cells.Item(44,5).Value = "  +9.35%  "

$ws.Cells.Item(45,4).Value = "0.9968"
$ws.Cells.Item(45,5).Value = "  +2.65%  "

$ws.Cells.Item(46,4).Value = "13.33"
$ws.Cells.Item(46,5).Value = "  +3.97%  "

$ws.Cells.Item(47,4).Value = "3.816"
$ws.Cells.Item(47,5).Value = "  +6.01%  "

$ws.Cells.Item(48,4).Value = "0.5993"
$ws.Cells.Item(48,5).Value = "  +6.65%  "

$ws.Cells.Item(49,4).Value = "127.86"
$ws.Cells.Item(49,5).Value = "  +2.92%  "

$ws.Cells.Item(50,4).Value = "2.045"

$ws.Cells.Item(51,4).Value = "0.07188"
$ws.Cells.Item(51,5).Value = "  +6.94%  "
